$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (row 37) was added to the data set. All the
# existing records that used to live in rows 37..80 shift down by one row
# (to rows 38..81); the worksheet's used-range grows from A1:R80 to A1:R81.
$ws.Rows.Item(37).Insert()

# Populate the newly inserted row 37 with the new record's values.
$ws.Range("A37").Value = 1
$ws.Range("B37").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C37").Value = "Arica y Parinacota"
$ws.Range("D37").Value = 44895
$ws.Range("E37").Value = 15
$ws.Range("F37").Value = 100112009
$ws.Range("G37").Value = "Acelga"
$ws.Range("H37").Value = "Sin especificar"
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 700
$ws.Range("K37").Value = 1700
$ws.Range("L37").Value = 2000
$ws.Range("M37").Value = 1850
$ws.Range("N37").Value = "`$/atado 2,5 a 3 kilos"
$ws.Range("O37").Value = "Región de Arica y Parinacota"
$ws.Range("P37").Value = 617
$ws.Range("Q37").Value = 3
$ws.Range("R37").Value = "Hortaliza"
